$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet is protected; lift protection so the cell values below can be
# written, then re-apply protection at the end.
$ws.Unprotect()

# Update the "as of" date in the confidentiality / disclosure note (A18).
$ws.Range("A18").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-13 for illustrative purposes only and are subject to change."

# Refresh the Weight (D) and Percent Change (E) figures for rows 2-15.
$ws.Range("D2").Value = 0.05574103990905283
$ws.Range("E2").Value = 0.0110204665807927

$ws.Range("D3").Value = 0.02343152702154676
$ws.Range("E3").Value = 0.006951340615690249

$ws.Range("D4").Value = 0.03052713456590678
$ws.Range("E4").Value = -0.001176009408075274

$ws.Range("D5").Value = 0.03295571644162561
$ws.Range("E5").Value = -0.009582215408202366

$ws.Range("D6").Value = 0.03932268294142285
$ws.Range("E6").Value = -0.01927368634611482

$ws.Range("D7").Value = 0.01917336471276518
$ws.Range("E7").Value = 0.01045388349514553

$ws.Range("D8").Value = 0.004130517162288837
$ws.Range("E8").Value = 0.003862868179623336

$ws.Range("D9").Value = 0.00679444477685368
$ws.Range("E9").Value = 0.02289628180039149

$ws.Range("D10").Value = 0.07211950776840383
$ws.Range("E10").Value = 0.004424778761061843

$ws.Range("D11").Value = 0.07219928598496181
$ws.Range("E11").Value = 0.004972375690607711

$ws.Range("D12").Value = 0.1442390155368077
$ws.Range("E12").Value = 0.001548672566371811

$ws.Range("D13").Value = 0.3848979836055764
$ws.Range("E13").Value = 0.003347427766032585

$ws.Range("D14").Value = 0.1144677795727876
$ws.Range("E14").Value = 0.002195389681668436

$ws.Range("D15").Value = 0.9999999999999999
$ws.Range("E15").Value = 0.002480763477532344

# Restore the sheet protection that was in place before the edits.
$ws.Protect()
